$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename "sheet10" -> "InputData", and add a brand new trailing sheet
#    named "sheet11" (mirrors workbook.xml sheet list change).
# ---------------------------------------------------------------------------
$inputData = $wb.Worksheets.Item("sheet10")
$inputData.Name = "InputData"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet11 = $wb.Worksheets.Add($null, $lastSheet)
$sheet11.Name = "sheet11"

# ---------------------------------------------------------------------------
# 2. InputData sheet (formerly sheet10): add A1 = "641001" (stored as text,
#    not auto-coerced to a number)
# ---------------------------------------------------------------------------
$inputData.Range("A1").NumberFormat = "@"
$inputData.Range("A1").Value = "641001"
$inputData.Range("A1").ClearFormats()

# ---------------------------------------------------------------------------
# 3. AssertData sheet: A2 "adidas" -> "ADIDAS", plus new rows 8-15
# ---------------------------------------------------------------------------
$assertData = $wb.Worksheets.Item("AssertData")
$assertData.Range("A2").Value = "ADIDAS"
$assertData.Range("A8").Value = "BADMINTON RACKETS"
$assertData.Range("A9").Value = "YONEX ARCSABER 2 FEEL"
$assertData.Range("A10").Value = "CART"
$assertData.Range("A11").Value = "ADDED TO CART"
$assertData.Range("A12").Value = "ADDED YONEX"
$assertData.Range("A13").Value = "ENTER ZIPCODE"
$assertData.Range("A14").Value = "shoes"
$assertData.Range("A15").Value = "required"

# ---------------------------------------------------------------------------
# 4. CreateAccountPage sheet: update A1/B1/A2 and add D7 (blank marker cell)
# ---------------------------------------------------------------------------
$createAccount = $wb.Worksheets.Item("CreateAccountPage")
$createAccount.Range("A1").Value = "wertyuilkjgfdsxcv"
$createAccount.Range("B1").Value = "fghjfds4@gmail.com"
$createAccount.Range("A2").Value = "dsfdfhgf"
$createAccount.Range("D7").Formula = "="""""
